## Adds 5 new registrant rows (20-24) to Sheet1, matching the data captured
## in the commit ("adding wp and meal option to the form", corrected mail
## sender id, cropped registration id, capitalised names, etc).
##
## Every column in this sheet is stored as literal text (Registration ID,
## phone numbers with leading zeros, year-only "dates", dd-like strings,
## etc.) so each cell's number format is forced to "@" (Text) before the
## value is poured in -- otherwise Excel/COM would happily "helpfully"
## coerce things like "2020" or "2002-01-21" into a real number/date and
## we would lose the leading zeros on the phone numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($sheet, $addr, $val) {
    $rng = $sheet.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
}

$newRows = @{
    20 = @{
        A = "5fda7191-7378-4204-a3ad-da7b84184725"
        B = "Agniva"
        C = "Bhattacharjee"
        D = "bvbnm,"
        E = "2020"
        F = "08420880979"
        G = "imagniva007@gmail.com"
        H = "2002-01-21"
        I = "hgfcgvhbn"
        J = "lkjkhghvhbn"
        K = "knjhghcvb"
        L = "mkjhgvhb"
        M = "mjhgfcgvb"
        N = "nbnvbvnm"
    }
    21 = @{
        A = "3c693800-cee7-4dfc-b8d6-d86c77512b5f"
        B = "AGNIVA"
        C = "BHATTACHARJEE"
        D = "biye to ei jonme hobe na"
        E = "1947"
        F = "08420880979"
        G = "bhattacharjee.agniva.jobs@gmail.com"
        H = "0089-05-04"
        I = "bekar jubok"
        J = "Google"
        K = "IAS OFFICER"
        L = "All India Bakchod"
        M = "nei kichu bhai"
        N = "dfgtgrfedsx"
    }
    22 = @{
        A = "BSS47afec85a74b"
        B = "AGNIVA"
        C = "BHATTACHARJEE"
        D = "nei amar"
        E = "2020"
        F = "08420880979"
        G = "bhattacharjee.agniva.jobs@gmail.com"
        H = "2002-01-21"
        I = "IT"
        J = "google"
        K = ""
        L = ""
        M = ""
        N = "sdsfdgfsvcs"
    }
    23 = @{
        A = "BSS - c51ccd599251"
        B = "Suchi"
        C = "BHATTACHARJEE"
        D = "nei amar"
        E = "2020"
        F = "08420880979"
        G = "bhattacharjee.agniva.jobs@gmail.com"
        H = "2002-01-21"
        I = "IT"
        J = "google"
        K = ""
        L = ""
        M = ""
        N = "sdsfdgfsvcs1223"
    }
    24 = @{
        A = "BSS - 9e614b382893"
        B = "Suchi"
        C = "BHATTACHARJEE"
        D = "nei amar"
        E = "2020"
        F = "08420880979"
        G = "bhattacharjee.agniva.jobs@gmail.com"
        H = "2002-01-21"
        I = "IT"
        J = "google"
        K = ""
        L = ""
        M = ""
        N = "sdsfdgfsvcs1223"
    }
}

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N")

foreach ($r in 20..24) {
    $rowData = $newRows[$r]
    foreach ($c in $cols) {
        Set-TextCell $ws "$c$r" $rowData[$c]
    }
}

## Row 24's phone number (F24) was entered as a formula in the source
## workbook (`=8420880979`) rather than a literal -- reproduce that.
$ws.Range("F24").NumberFormat = "@"
$ws.Range("F24").Formula = "=8420880979"

## Keep the worksheet's "numbers stored as text" warning suppressed over
## the full, now-larger, data range (A1:N24), matching how the original
## sheet suppressed it over A1:N19.
try {
    $ws.Range("A1:N24").Errors.Item(9).Ignore = $true
} catch {
    # Older/limited hosts may not expose the NumberAsText error-checking
    # item; harmless to skip since it only toggles a UI warning glyph.
}
